# Add files via upload
# Renumber the "ID" column (column B) for rows 29-33 on sheet SRC_CZ:
# a new row was inserted logically at row 29 (ID 27), so the existing
# IDs 27..30 (previously in rows 29..32) are shifted down by one and a
# new ID (27) is written into the now-empty B29 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(29, 2).Value = 27
$ws.Cells.Item(30, 2).Value = 28
$ws.Cells.Item(31, 2).Value = 29
$ws.Cells.Item(32, 2).Value = 30
$ws.Cells.Item(33, 2).Value = 31

# Reflect the selection state recorded in the workbook (B26:B38, active cell B26)
$ws.Range("B26:B38").Select()
